# Add basic design, table and buttons to DriverOverview Component
#
# The plan table on Sheet1 gets its last task row re-purposed as a
# "Map preview" task, and a new row is appended for a
# "Driver preview module" task.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 previously held the "Create Driver Tracking module" task;
# re-label it as the Map preview task.
$ws.Range("B12").Value = "Create Map Previev Component"

# Append a new row 13 for the driver preview module task.
$ws.Range("A13").Value = "#100011"
$ws.Range("B13").Value = "Create Driver preview module"

# Restore the (somewhat stale) selection state recorded in the workbook.
$ws.Range("B23").Select()
